# Auto-generated script to apply cell-value updates per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1685.9143
$ws.Range("J17").Value = 1854.1154
$ws.Range("L17").Value = 5562.3462
$ws.Range("N17").Value = -5898.3462
$ws.Range("H29").Value = 913.8570999999999
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15562
$ws.Range("H62").Value = 16151.44
$ws.Range("I62").Value = 15869.444
$ws.Range("K62").Value = 15869.444
$ws.Range("M62").Value = -15245.444
$ws.Range("H65").Value = 16151.44
$ws.Range("I65").Value = 15869.444
$ws.Range("K65").Value = 79347.22
$ws.Range("M65").Value = -76227.22
$ws.Range("H70").Value = 7434.4116
$ws.Range("I70").Value = 1648.875
$ws.Range("J70").Value = 100003
$ws.Range("K70").Value = 4946.625
$ws.Range("L70").Value = 300009
$ws.Range("M70").Value = -4676.625
$ws.Range("N70").Value = -300549
$ws.Range("H73").Value = 7434.4116
$ws.Range("I73").Value = 1648.875
$ws.Range("J73").Value = 100003
$ws.Range("K73").Value = 4946.625
$ws.Range("L73").Value = 300009
$ws.Range("M73").Value = -4010.625
$ws.Range("N73").Value = -301881
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("H112").Value = 1922.1852
$ws.Range("J112").Value = 1995.8334
$ws.Range("L112").Value = 5987.5002
$ws.Range("N112").Value = -8203.5002
$ws.Range("H116").Value = 42278.375
$ws.Range("J116").Value = 3862.25
$ws.Range("L116").Value = 3862.25
$ws.Range("N116").Value = -10746.25
$ws.Range("H138").Value = 3692.5244
$ws.Range("J138").Value = 4047.209
$ws.Range("L138").Value = 12141.627
$ws.Range("N138").Value = -22421.627

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4169741.8
$ws.Range("I74").Value = 6252164.5
$ws.Range("K74").Value = 6252164.5
$ws.Range("M74").Value = -6251290.5
$ws.Range("H77").Value = 4169741.8
$ws.Range("I77").Value = 6252164.5
$ws.Range("K77").Value = 31260822.5
$ws.Range("M77").Value = -31256454.5
$ws.Range("H102").Value = 1784.8918
$ws.Range("I102").Value = 1872.1212
$ws.Range("K102").Value = 1872.1212
$ws.Range("M102").Value = -250.1212
$ws.Range("H110").Value = 2246.7
$ws.Range("I110").Value = 2290.4119
$ws.Range("J110").Value = 1999
$ws.Range("K110").Value = 2290.4119
$ws.Range("L110").Value = 1999
$ws.Range("M110").Value = -245.4119000000001
$ws.Range("N110").Value = -6089

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 12241.667
$ws.Range("I26").Value = 12241.667
$ws.Range("K26").Value = 12241.667
$ws.Range("M26").Value = -11949.667
$ws.Range("H94").Value = 1261.2
$ws.Range("I94").Value = 1436.5
$ws.Range("J94").Value = 998.25
$ws.Range("K94").Value = 1436.5
$ws.Range("L94").Value = 998.25
$ws.Range("M94").Value = -985.5
$ws.Range("N94").Value = -1900.25
$ws.Range("H105").Value = 4058.889
$ws.Range("I105").Value = 3816.375
$ws.Range("K105").Value = 3816.375
$ws.Range("M105").Value = -2069.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13342.071
$ws.Range("I31").Value = 2499.5
$ws.Range("K31").Value = 2499.5
$ws.Range("M31").Value = -2204.5
$ws.Range("H34").Value = 13342.071
$ws.Range("I34").Value = 2499.5
$ws.Range("K34").Value = 2499.5
$ws.Range("M34").Value = -2297.5
$ws.Range("H134").Value = 4509.7173
$ws.Range("I134").Value = 1591.907
$ws.Range("K134").Value = 4775.721
$ws.Range("M134").Value = -2240.721

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 618.6667
$ws.Range("I8").Value = 618.6667
$ws.Range("K8").Value = 1856.0001
$ws.Range("M8").Value = -1717.0001
$ws.Range("H56").Value = 9598
$ws.Range("I56").Value = 9598
$ws.Range("K56").Value = 9598
$ws.Range("M56").Value = -9068
$ws.Range("H62").Value = 6497
$ws.Range("J62").Value = 6497
$ws.Range("L62").Value = 19491
$ws.Range("N62").Value = -20863
$ws.Range("H65").Value = 6497
$ws.Range("J65").Value = 6497
$ws.Range("L65").Value = 58473
$ws.Range("N65").Value = -65337
$ws.Range("H82").Value = 8656.333000000001
$ws.Range("I82").Value = 8000
$ws.Range("K82").Value = 24000
$ws.Range("M82").Value = -23594
$ws.Range("H85").Value = 8656.333000000001
$ws.Range("I85").Value = 8000
$ws.Range("K85").Value = 24000
$ws.Range("M85").Value = -22596
$ws.Range("H86").Value = 1102.4286
$ws.Range("J86").Value = 1104.25
$ws.Range("L86").Value = 3312.75
$ws.Range("N86").Value = -5684.75
$ws.Range("H89").Value = 1102.4286
$ws.Range("J89").Value = 1104.25
$ws.Range("L89").Value = 9938.25
$ws.Range("N89").Value = -21794.25
$ws.Range("H107").Value = 966.3333
$ws.Range("I107").Value = 749.5
$ws.Range("J107").Value = 1400
$ws.Range("K107").Value = 2248.5
$ws.Range("L107").Value = 4200
$ws.Range("M107").Value = -328.5
$ws.Range("N107").Value = -8040
$ws.Range("H113").Value = 1497.7142
$ws.Range("J113").Value = 1516.95
$ws.Range("L113").Value = 4550.85
$ws.Range("N113").Value = -8890.85
$ws.Range("H129").Value = 3392.5833
$ws.Range("J129").Value = 3628.3635
$ws.Range("L129").Value = 10885.0905
$ws.Range("N129").Value = -20885.0905
$ws.Range("H131").Value = 7327.1284
$ws.Range("I131").Value = 983.75
$ws.Range("J131").Value = 10146.407
$ws.Range("K131").Value = 2951.25
$ws.Range("L131").Value = 30439.221
$ws.Range("M131").Value = 2088.75
$ws.Range("N131").Value = -40519.221
$ws.Range("H132").Value = 3849.6667
$ws.Range("I132").Value = 2834.2
$ws.Range("J132").Value = 5119
$ws.Range("K132").Value = 25507.8
$ws.Range("L132").Value = 46071
$ws.Range("M132").Value = -22977.8
$ws.Range("N132").Value = -51131
$ws.Range("H137").Value = 4328.8335
$ws.Range("I137").Value = 3995.75
$ws.Range("K137").Value = 11987.25
$ws.Range("M137").Value = -6887.25
$ws.Range("H141").Value = 2432.125
$ws.Range("I141").Value = 1922.4286
$ws.Range("K141").Value = 5767.2858
$ws.Range("M141").Value = -587.2857999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 19087
$ws.Range("J32").Value = 19087
$ws.Range("L32").Value = 19087
$ws.Range("N32").Value = -19679
$ws.Range("H80").Value = 5100.077
$ws.Range("I80").Value = 5350.25
$ws.Range("K80").Value = 5350.25
$ws.Range("M80").Value = -4352.25
$ws.Range("H83").Value = 5100.077
$ws.Range("I83").Value = 5350.25
$ws.Range("K83").Value = 26751.25
$ws.Range("M83").Value = -21759.25
$ws.Range("H102").Value = 2316.4482
$ws.Range("I102").Value = 2174.0417
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2174.0417
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -552.0417000000002
$ws.Range("N102").Value = -6244
$ws.Range("H126").Value = 836327.25
$ws.Range("I126").Value = 1391612.8
$ws.Range("J126").Value = 3399
$ws.Range("K126").Value = 4174838.4
$ws.Range("L126").Value = 10197
$ws.Range("M126").Value = -4172368.4
$ws.Range("N126").Value = -15137

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 30870
$ws.Range("I22").Value = 92070.27
$ws.Range("J22").Value = 2819.875
$ws.Range("K22").Value = 92070.27
$ws.Range("L22").Value = 2819.875
$ws.Range("M22").Value = -91775.27
$ws.Range("N22").Value = -3409.875
$ws.Range("H27").Value = 30870
$ws.Range("I27").Value = 92070.27
$ws.Range("J27").Value = 2819.875
$ws.Range("K27").Value = 92070.27
$ws.Range("L27").Value = 2819.875
$ws.Range("M27").Value = -91963.27
$ws.Range("N27").Value = -3033.875
$ws.Range("H40").Value = 1504.6428
$ws.Range("I40").Value = 1504.6428
$ws.Range("K40").Value = 1504.6428
$ws.Range("M40").Value = -1368.6428
$ws.Range("H68").Value = 9997.200000000001
$ws.Range("I68").Value = 9997.200000000001
$ws.Range("K68").Value = 9997.200000000001
$ws.Range("M68").Value = -9248.200000000001
$ws.Range("H71").Value = 9997.200000000001
$ws.Range("I71").Value = 9997.200000000001
$ws.Range("K71").Value = 49986
$ws.Range("M71").Value = -46242
$ws.Range("H100").Value = 6743.9614
$ws.Range("I100").Value = 1569.1538
$ws.Range("K100").Value = 1569.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1846.4445
$ws.Range("J100").Value = 1139.2
$ws.Range("L100").Value = 2278.4
$ws.Range("N100").Value = -3360.4
$ws.Range("H122").Value = 3816.9473
$ws.Range("I122").Value = 3434.8667
$ws.Range("K122").Value = 10304.6001
$ws.Range("M122").Value = -7854.6001
